# Fruta / hortaliza, semanal
# Weekly refresh: rows 3-9 (the per-record detail rows, row 1 = header,
# row 2 = first unchanged record) get reshuffled with updated data for
# Variedad/Calidad/Volumen/Precio*/Origen/Precio-$-Kg, keyed by the row's
# Fecha (column D). Row 2 is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to write into rows 3..9, column-by-column.
# Columns: D=Fecha, K=Variedad, L=Calidad, M=Volumen, N=Precio minimo,
#          O=Precio maximo, P=Precio promedio ponderado, R=Origen,
#          S=Precio $/Kg
$rows = @(
    @{ Row = 3; D = 44342; K = "Mankaki"; L = "Primera"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 },
    @{ Row = 4; D = 44301; K = "Hachiya"; L = "Segunda"; M = 250; N = 20000; O = 21000; P = 20500; R = "Región de O'Higgins"; S = 1139 },
    @{ Row = 5; D = 44699; K = "Mankaki"; L = "Primera"; M = 250; N = 29000; O = 30000; P = 29500; R = "Región de O'Higgins"; S = 1639 },
    @{ Row = 6; D = 44305; K = "Mankaki"; L = "Segunda"; M = 250; N = 24000; O = 25000; P = 24500; R = "Región de O'Higgins"; S = 1361 },
    @{ Row = 7; D = 44313; K = "Mankaki"; L = "Primera"; M = 270; N = 21000; O = 22000; P = 21500; R = "Región de O'Higgins"; S = 1194 },
    @{ Row = 8; D = 45043; K = "Fuyu";    L = "Primera"; M = 300; N = 25000; O = 26000; P = 25500; R = "Región de O'Higgins"; S = 1417 },
    @{ Row = 9; D = 45071; K = "Fuyu";    L = "Segunda"; M = 110; N = 23000; O = 24000; P = 23455; R = "Región Metropolitana"; S = 1303 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($row, 11).Value = $r.K   # K - Variedad
    $ws.Cells.Item($row, 12).Value = $r.L   # L - Calidad
    $ws.Cells.Item($row, 13).Value = $r.M   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $r.N   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $r.O   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $r.P   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 18).Value = $r.R   # R - Origen
    $ws.Cells.Item($row, 19).Value = $r.S   # S - Precio $/Kg
}
